$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value = 8626548
$ws.Range("I62").Value = 12507415
$ws.Range("K62").Value = 12507415
$ws.Range("M62").Value = -12506791
# row 65
$ws.Range("H65").Value = 8626548
$ws.Range("I65").Value = 12507415
$ws.Range("K65").Value = 62537075
$ws.Range("M65").Value = -62533955
# row 132
$ws.Range("H132").Value = 11635672
$ws.Range("I132").Value = 17249356
$ws.Range("J132").Value = 7328.357
$ws.Range("K132").Value = 51748068
$ws.Range("L132").Value = 21985.071
$ws.Range("M132").Value = -51745538
$ws.Range("N132").Value = -27045.071
# row 137
$ws.Range("H137").Value = 1288.4667
$ws.Range("I137").Value = 1188.6552
$ws.Range("J137").Value = 1469.375
$ws.Range("K137").Value = 3565.9656
$ws.Range("L137").Value = 4408.125
$ws.Range("M137").Value = -1015.9656
$ws.Range("N137").Value = -9508.125

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 59913
$ws.Range("J45").Value = 1158.1666
$ws.Range("L45").Value = 1158.1666
$ws.Range("N45").Value = -1912.1666
# row 125
$ws.Range("H125").Value = 61983
$ws.Range("J125").Value = 61983
$ws.Range("L125").Value = 61983
$ws.Range("N125").Value = -71823

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 2117073.5
$ws.Range("I86").Value = 3783.3333
$ws.Range("J86").Value = 4653022
$ws.Range("K86").Value = 3783.3333
$ws.Range("L86").Value = 4653022
$ws.Range("M86").Value = -2660.3333
$ws.Range("N86").Value = -4655268
# row 89
$ws.Range("H89").Value = 2117073.5
$ws.Range("I89").Value = 3783.3333
$ws.Range("J89").Value = 4653022
$ws.Range("K89").Value = 18916.6665
$ws.Range("L89").Value = 23265110
$ws.Range("M89").Value = -13300.6665
$ws.Range("N89").Value = -23276342
# row 105
$ws.Range("H105").Value = 90910990
$ws.Range("I105").Value = 1888.75
$ws.Range("K105").Value = 1888.75
$ws.Range("M105").Value = -141.75
# row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# row 132
$ws.Range("H132").Value = 44892.727
$ws.Range("J132").Value = 44892.727
$ws.Range("L132").Value = 44892.727
$ws.Range("N132").Value = -55012.727

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1469.561
$ws.Range("I31").Value = 1238.2858
$ws.Range("J31").Value = 1589.4814
$ws.Range("K31").Value = 1238.2858
$ws.Range("L31").Value = 1589.4814
$ws.Range("M31").Value = -943.2858000000001
$ws.Range("N31").Value = -2179.4814
# row 34
$ws.Range("H34").Value = 1469.561
$ws.Range("I34").Value = 1238.2858
$ws.Range("J34").Value = 1589.4814
$ws.Range("K34").Value = 1238.2858
$ws.Range("L34").Value = 1589.4814
$ws.Range("M34").Value = -1036.2858
$ws.Range("N34").Value = -1993.4814
# row 96
$ws.Range("H96").Value = 17464.5
$ws.Range("J96").Value = 17464.5
$ws.Range("L96").Value = 17464.5
$ws.Range("N96").Value = -22956.5
# row 134
$ws.Range("H134").Value = 1211
$ws.Range("I134").Value = 1201.4286
$ws.Range("J134").Value = 1345
$ws.Range("K134").Value = 3604.2858
$ws.Range("L134").Value = 4035
$ws.Range("M134").Value = -1069.2858
$ws.Range("N134").Value = -9105

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 103
$ws.Range("H103").Value = 3623.5454
$ws.Range("J103").Value = 4376
$ws.Range("L103").Value = 13128
$ws.Range("N103").Value = -14886
# row 113
$ws.Range("H113").Value = 27193410
$ws.Range("I113").Value = 20833632
$ws.Range("J113").Value = 28889350
$ws.Range("K113").Value = 62500896
$ws.Range("L113").Value = 86668050
$ws.Range("M113").Value = -62498726
$ws.Range("N113").Value = -86672390
# row 131
$ws.Range("H131").Value = 931.3196
$ws.Range("J131").Value = 973.86365
$ws.Range("L131").Value = 2921.59095
$ws.Range("N131").Value = -13001.59095

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 5647.0586
$ws.Range("I70").Value = 4400
$ws.Range("J70").Value = 7428.5713
$ws.Range("K70").Value = 4400
$ws.Range("L70").Value = 7428.5713
$ws.Range("M70").Value = -4130
$ws.Range("N70").Value = -7968.5713
# row 73
$ws.Range("H73").Value = 5647.0586
$ws.Range("I73").Value = 4400
$ws.Range("J73").Value = 7428.5713
$ws.Range("K73").Value = 4400
$ws.Range("L73").Value = 7428.5713
$ws.Range("M73").Value = -3464
$ws.Range("N73").Value = -9300.5713
# row 132
$ws.Range("H132").Value = 8676.16
$ws.Range("I132").Value = 5464.476
$ws.Range("K132").Value = 16393.428
$ws.Range("M132").Value = -13863.428
# row 135
$ws.Range("H135").Value = 58797.777
$ws.Range("J135").Value = 58797.777
$ws.Range("L135").Value = 58797.777
$ws.Range("N135").Value = -68937.777

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 20005904
$ws.Range("I132").Value = 28573140
$ws.Range("J132").Value = 15689.733
$ws.Range("K132").Value = 85719420
$ws.Range("L132").Value = 47069.199
$ws.Range("M132").Value = -85716890
$ws.Range("N132").Value = -52129.199

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 46
$ws.Range("H46").Value = 24500
$ws.Range("J46").Value = 24500
$ws.Range("L46").Value = 24500
$ws.Range("N46").Value = -24962
# row 81
$ws.Range("H81").Value = 1068.2
$ws.Range("I81").Value = 1068.2
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2136.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1075.4
$ws.Range("N81").ClearContents()
# row 84
$ws.Range("H84").Value = 1068.2
$ws.Range("I84").Value = 1068.2
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10682
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -5378
$ws.Range("N84").ClearContents()
# row 132
$ws.Range("H132").Value = 17846.516
$ws.Range("I132").Value = 23663.844
$ws.Range("J132").Value = 6464.7827
$ws.Range("K132").Value = 70991.53200000001
$ws.Range("L132").Value = 19394.3481
$ws.Range("M132").Value = -68461.53200000001
$ws.Range("N132").Value = -24454.3481
# row 134
$ws.Range("H134").Value = 24500
$ws.Range("J134").Value = 24500
$ws.Range("L134").Value = 73500
$ws.Range("N134").Value = -78570
# row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# row 140
$ws.Range("H140").Value = 37929.465
$ws.Range("J140").Value = 37929.465
$ws.Range("L140").Value = 37929.465
$ws.Range("N140").Value = -48289.465
